$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "HP" brand row (row 5) was moved down to the bottom of the list
# (row 31). Delete the row at its old position; this shifts every
# subsequent row up by one. Then write "HP" into both columns of the
# now-last row.
$ws.Rows("5").Delete()
$ws.Range("A31").Value = "HP"
$ws.Range("B31").Value = "HP"
